$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chris")

$ws.Range("A19").Value = "x"
$ws.Range("A19").Font.Name = "Calibri"
$ws.Range("A19").Interior.Pattern = -4142
$ws.Range("A19").HorizontalAlignment = -4152

$ws.Range("B19").Font.Name = "Calibri"
$ws.Range("B19").Interior.Pattern = -4142
$ws.Range("B19").HorizontalAlignment = -4108
